$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value2 = 0.6352250347264584
$ws.Range("J2").Value2 = 0.6352250347264584
$ws.Range("M2").Value2 = 2.598166333333333
$ws.Range("N2").Value2 = 7.794499
$ws.Range("O2").Value2 = 0.3466013321552429
$ws.Range("P2").Value2 = 0.3466013321552429
$ws.Range("Q2").Value2 = 20.86880109040222
$ws.Range("R2").Value2 = 187.81920981362
$ws.Range("S2").Value2 = 0.2201698432545509
$ws.Range("T2").Value2 = 0.2201698432545509
$ws.Range("I3").Value2 = 0.6352250347264584
$ws.Range("J3").Value2 = 0.6352250347264584
$ws.Range("M3").Value2 = 4.333403333333333
$ws.Range("O3").Value2 = 0.5780859172985858
$ws.Range("P3").Value2 = 0.5780859172985858
$ws.Range("S3").Value2 = 0.3672146468908707
$ws.Range("T3").Value2 = 0.3672146468908707
$ws.Range("I4").Value2 = 0.6352250347264584
$ws.Range("J4").Value2 = 0.6352250347264584
$ws.Range("M4").Value2 = 0.4692043333333333
$ws.Range("N4").Value2 = 1.407613
$ws.Range("O4").Value2 = 0.06259293136852516
$ws.Range("P4").Value2 = 0.06259293136852516
$ws.Range("Q4").Value2 = 3.768708637882222
$ws.Range("R4").Value2 = 33.91837774094
$ws.Range("S4").Value2 = 0.03976059700220222
$ws.Range("T4").Value2 = 0.03976059700220222
$ws.Range("I5").Value2 = 0.6352250347264584
$ws.Range("J5").Value2 = 0.6352250347264584
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 0.09534933333333333
$ws.Range("N5").Value2 = 0.286048
$ws.Range("O5").Value2 = 0.01271981917764605
$ws.Range("P5").Value2 = 0.01271981917764604
$ws.Range("Q5").Value2 = 0.7658579229155555
$ws.Range("R5").Value2 = 6.892721306239999
$ws.Range("S5").Value2 = 0.008079947578834481
$ws.Range("T5").Value2 = 0.008079947578834481
$ws.Range("G6").Value2 = 4.495264666666666
$ws.Range("H6").Value2 = 13.485794
$ws.Range("I6").Value2 = 0.3555104111888949
$ws.Range("J6").Value2 = 0.3555104111888949
$ws.Range("M6").Value2 = 2.598166333333333
$ws.Range("N6").Value2 = 7.794499
$ws.Range("O6").Value2 = 0.3466013321552429
$ws.Range("P6").Value2 = 0.3466013321552429
$ws.Range("Q6").Value2 = 11.67944531635622
$ws.Range("R6").Value2 = 105.115007847206
$ws.Range("S6").Value2 = 0.1232203821131292
$ws.Range("T6").Value2 = 0.1232203821131292
$ws.Range("G7").Value2 = 4.495264666666666
$ws.Range("H7").Value2 = 13.485794
$ws.Range("I7").Value2 = 0.3555104111888949
$ws.Range("J7").Value2 = 0.3555104111888949
$ws.Range("M7").Value2 = 4.333403333333333
$ws.Range("O7").Value2 = 0.5780859172985858
$ws.Range("P7").Value2 = 0.5780859172985858
$ws.Range("Q7").Value2 = 19.47979489074889
$ws.Range("R7").Value2 = 175.31815401674
$ws.Range("S7").Value2 = 0.2055155621613297
$ws.Range("T7").Value2 = 0.2055155621613297
$ws.Range("G8").Value2 = 4.495264666666666
$ws.Range("H8").Value2 = 13.485794
$ws.Range("I8").Value2 = 0.3555104111888949
$ws.Range("J8").Value2 = 0.3555104111888949
$ws.Range("M8").Value2 = 0.4692043333333333
$ws.Range("N8").Value2 = 1.407613
$ws.Range("O8").Value2 = 0.06259293136852516
$ws.Range("P8").Value2 = 0.06259293136852516
$ws.Range("Q8").Value2 = 2.109197661080222
$ws.Range("R8").Value2 = 18.982778949722
$ws.Range("S8").Value2 = 0.02225243876834266
$ws.Range("T8").Value2 = 0.02225243876834266
$ws.Range("G9").Value2 = 4.495264666666666
$ws.Range("H9").Value2 = 13.485794
$ws.Range("I9").Value2 = 0.3555104111888949
$ws.Range("J9").Value2 = 0.3555104111888949
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 0.6666666666666666
$ws.Range("M9").Value2 = 0.09534933333333333
$ws.Range("N9").Value2 = 0.286048
$ws.Range("O9").Value2 = 0.01271981917764605
$ws.Range("P9").Value2 = 0.01271981917764604
$ws.Range("Q9").Value2 = 0.4286204891235555
$ws.Range("R9").Value2 = 3.857584402112
$ws.Range("S9").Value2 = 0.004522028146093337
$ws.Range("T9").Value2 = 0.004522028146093336
$ws.Range("E10").Value2 = 2
$ws.Range("F10").Value2 = 0.6666666666666666
$ws.Range("G10").Value2 = 0.117146
$ws.Range("H10").Value2 = 0.351438
$ws.Range("I10").Value2 = 0.009264554084646619
$ws.Range("J10").Value2 = 0.009264554084646619
$ws.Range("M10").Value2 = 2.598166333333333
$ws.Range("N10").Value2 = 7.794499
$ws.Range("O10").Value2 = 0.3466013321552429
$ws.Range("P10").Value2 = 0.3466013321552429
$ws.Range("Q10").Value2 = 0.3043647932846667
$ws.Range("R10").Value2 = 2.739283139562
$ws.Range("S10").Value2 = 0.003211106787562815
$ws.Range("T10").Value2 = 0.003211106787562815
$ws.Range("E11").Value2 = 2
$ws.Range("F11").Value2 = 0.6666666666666666
$ws.Range("G11").Value2 = 0.117146
$ws.Range("H11").Value2 = 0.351438
$ws.Range("I11").Value2 = 0.009264554084646619
$ws.Range("J11").Value2 = 0.009264554084646619
$ws.Range("M11").Value2 = 4.333403333333333
$ws.Range("O11").Value2 = 0.5780859172985858
$ws.Range("P11").Value2 = 0.5780859172985858
$ws.Range("Q11").Value2 = 0.5076408668866667
$ws.Range("R11").Value2 = 4.56876780198
$ws.Range("S11").Value2 = 0.005355708246385301
$ws.Range("T11").Value2 = 0.005355708246385301
$ws.Range("E12").Value2 = 2
$ws.Range("F12").Value2 = 0.6666666666666666
$ws.Range("G12").Value2 = 0.117146
$ws.Range("H12").Value2 = 0.351438
$ws.Range("I12").Value2 = 0.009264554084646619
$ws.Range("J12").Value2 = 0.009264554084646619
$ws.Range("M12").Value2 = 0.4692043333333333
$ws.Range("N12").Value2 = 1.407613
$ws.Range("O12").Value2 = 0.06259293136852516
$ws.Range("P12").Value2 = 0.06259293136852516
$ws.Range("Q12").Value2 = 0.05496541083266667
$ws.Range("R12").Value2 = 0.494688697494
$ws.Range("S12").Value2 = 0.0005798955979802753
$ws.Range("T12").Value2 = 0.0005798955979802753
$ws.Range("E13").Value2 = 2
$ws.Range("F13").Value2 = 0.6666666666666666
$ws.Range("G13").Value2 = 0.117146
$ws.Range("H13").Value2 = 0.351438
$ws.Range("I13").Value2 = 0.009264554084646619
$ws.Range("J13").Value2 = 0.009264554084646619
$ws.Range("K13").Value2 = 2
$ws.Range("L13").Value2 = 0.6666666666666666
$ws.Range("M13").Value2 = 0.09534933333333333
$ws.Range("N13").Value2 = 0.286048
$ws.Range("O13").Value2 = 0.01271981917764605
$ws.Range("P13").Value2 = 0.01271981917764604
$ws.Range("Q13").Value2 = 0.01116979300266667
$ws.Range("R13").Value2 = 0.100528137024
$ws.Range("S13").Value2 = 0.0001178434527182271
$ws.Range("T13").Value2 = 0.0001178434527182271
